$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.291.02"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.858.69"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'0.7030"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").Value = "'238.36"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.07903"
$ws.Range("E8").Value = "  +3.23%  "
$ws.Range("D9").Value = "'0.3044"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "'24.46"
$ws.Range("E10").Value = "  +5.45%  "
$ws.Range("D11").Value = "'0.08179"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7234"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.216"
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'89.64"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.698.19"
$ws.Range("E15").Value = "  -8.66%  "
$ws.Range("D16").Value = "29.303.33"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "'5.818"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").Value = "'0.000007818"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").Value = "'13.25"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "'238.59"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "'7.555"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "1.871.48"
$ws.Range("E24").Value = "  -11.19%  "
$ws.Range("D25").Value = "'162.13"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").Value = "'8.904"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("E27").Value = "  -2.62%  "
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("D29").Value = "'1.916"
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("D30").Value = "'1.384"
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("D33").Value = "'4.057"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("D34").Value = "'0.05182"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("D36").Value = "'0.7114"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "'0.9961"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "'2.690"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("D41").Value = "1.158.47"
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("D42").Value = "'0.9212"
$ws.Range("E42").Value = "  -1.12%  "
$ws.Range("D43").Value = "'5.918"
$ws.Range("E43").Value = "  +1.12%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'71.07"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.4257"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "'100.97"
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("D48").Value = "'0.5317"
$ws.Range("E48").Value = "  -2.90%  "
$ws.Range("D49").Value = "'1.756"
$ws.Range("E49").Value = "  -1.98%  "
$ws.Range("D50").Value = "'9.194"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").Value = "'6.999"
$ws.Range("E51").Value = "  +0.76%  "
